$wb = $excel.ActiveWorkbook

# 1. Collection sheet: update version string in B2
$wsCollection = $wb.Worksheets.Item("Collection")
$wsCollection.Range("B2").Value = "Sample dataset for DataCrate v1.0"

# 2. Equipment sheet: remove the ImageMagick row (row 4)
$wsEquipment = $wb.Worksheets.Item("Equipment")
$wsEquipment.Rows.Item(4).Delete()

# 3. Actions sheet: add Agent value to H3
$wsActions = $wb.Worksheets.Item("Actions")
$wsActions.Range("H3").Value = "Peter Sefton"
